# Updated symbol list on Fri Dec 23 14:56:31 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores values as text (e.g. "244.75"), not as
# numbers. Excel auto-converts a numeric-looking string assigned via
# .Value into a real number, so those cells are first force-formatted as
# Text ("@") to preserve the original text semantics of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric - keep them stored as text like the
# rest of column D by pre-setting a Text number format.
$textFormatCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "244.67"
$ws.Range("D3").Value = "21.94"
$ws.Range("D4").Value = "5.394"
$ws.Range("D5").Value = "0.05848"
$ws.Range("D6").Value = "3.392"
$ws.Range("D7").Value = "6.359"
$ws.Range("D8").Value = "0.8151"
$ws.Range("D9").Value = "1.020"
$ws.Range("D10").Value = "0.1421"
$ws.Range("D11").Value = "0.03704"
$ws.Range("D12").Value = "0.07432"
$ws.Range("D13").Value = "0.03047"
$ws.Range("D14").Value = "4.200"
$ws.Range("D15").Value = "0.09394"
$ws.Range("D16").Value = "0.001600"
$ws.Range("D17").Value = "0.04805"
$ws.Range("D18").Value = "0.0005897"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "0.006067"

$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "0.001003"
$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "0.004083"
$ws.Range("E21").Value = "20HotbitTokenHTB"

$ws.Range("D22").Value = "0.0001502"
$ws.Range("D23").Value = "3.693"
$ws.Range("D24").Value = "2.224"
$ws.Range("D25").Value = "0.3232"
$ws.Range("D26").Value = "0.1296"
$ws.Range("D27").Value = "0.0002904"
$ws.Range("E27").Value = "26UpBotsUBXTWorstin24h"

$ws.Range("D40").Value = "0.03847"
$ws.Range("D41").Value = "0.006368"
$ws.Range("D42").Value = "0.1071"
$ws.Range("D43").Value = "0.003004"
$ws.Range("D44").Value = "0.006244"
$ws.Range("D45").Value = "0.00005629"
$ws.Range("D47").Value = "0.8110"
$ws.Range("D48").Value = "0.1425"
$ws.Range("D49").Value = "0.00002103"
